$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the style of an existing header cell (AC1) to the new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in Wins/Losses/Ties values for each data row (2 through 43)
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 30).Value = 100   # AD
    $ws.Cells.Item($r, 31).Value = 62    # AE
    $ws.Cells.Item($r, 32).Value = 0     # AF
}
